$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the cell before editing it, mirroring the on-screen navigation
# that would normally happen when a user scrolls to column U and clicks U1.
$ws.Range("U1").Select()

# Change the header text from "Approver" to "Requester"
$ws.Range("U1").Value = "Requester"

# Reflect the resulting view state: scrolled so column J is the left-most
# visible column, with U2 as the active selection.
$ws.Range("U2").Select()
$excel.ActiveWindow.ScrollColumn = 10
